$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers become numeric indices 0..5 (keep existing bold/centered style)
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5

# Column A becomes very narrow
$ws.Range("A1").EntireColumn.ColumnWidth = 0.15

# New row 5 of data
$ws.Range("A5").Value = "e350f4f4-91d9-47d2-8ada-c2c35e6b7926"
$ws.Range("B5").Value = 45660.77291162037
$ws.Range("C5").Value = "Cj"
$ws.Range("D5").Value = 45660.77291162037
$ws.Range("E5").Value = "11d8663c-2180-4d4c-94e6-d387090cbc0c"
$ws.Range("F5").Value = "d4e541ab-52c8-4d0d-a28f-cd5f3856aec6"

# Match date/time number format used by the other rows' timestamp columns
$ws.Range("B5").NumberFormat = $ws.Range("B4").NumberFormat
$ws.Range("D5").NumberFormat = $ws.Range("D4").NumberFormat

Write-Output "done"
